$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the JSON payload in B3 to include the new "times" field
$ws.Range("B3").Value = '{"userName":"hugang","caseId":"hugangのテスト","times":"5"}'

# E2 used to hold the literal number 5; it now references the ${times} placeholder
$ws.Range("E2").Value = '${times}'

# Update the if-condition in D2 to also check the new "times" variable
$ws.Range("D2").Value = "userName == 'hugang' && parseInt(times) == 5"

# Adjust column widths to fit the new, longer content
$ws.Columns.Item(2).ColumnWidth = 59.75
$ws.Columns.Item(4).ColumnWidth = 46.5
$ws.Columns.Item(5).ColumnWidth = 8.875

# Move the active selection
$ws.Range("D6").Select()
